$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The bitacora (log) table keeps one "entry" every two rows: the first row of
# the pair carries the date + comment text, the second row is a thick-bottom
# spacer row. Rows 19:20 hold the most recent existing entry; rows 21:22 are
# the next (still-empty) template pair. Copy the formatting of the existing
# entry down onto the new pair so the new row inherits the same fonts,
# borders, number format and merged cells, then fill in the new log entry.
$ws.Range("C19:E20").Copy()
$ws.Range("C21:E22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C21").Value2 = "Views"
$ws.Range("D21").Value2 = 44737
$ws.Range("E21").Value2 = "Creacion de Views"
$ws.Rows.Item(21).RowHeight = 15

# Match the author's recorded view state: scrolled a bit further down the
# sheet with cell H22 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H22").Select()
